# Daily attendance processing - 2025-12-24 22:30:33
# Normalizes the "Recorded By" column (G) so that any "System"/"system"
# entries are moved to the front of the comma-separated list, while the
# relative order of the remaining entries (and of the System-like tokens
# themselves) is reversed along with the rest of the list.
#
# Rows whose "Recorded By" value has no "System" token, or only a single
# token (nothing to reorder), are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }

    $parts = $val -split ','
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -lt 2) { continue }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.ToLower() -eq 'system') { $hasSystem = $true }
    }
    if (-not $hasSystem) { continue }

    [array]::Reverse($parts)
    $newVal = [string]::Join(', ', $parts)

    $cell.Value2 = $newVal
}
